$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# The row for Caso 2620 (DORREGO 2398) was removed from the source data.
# Deleting the entire row shifts every subsequent row up by one and
# reduces the used range from A1:N40 to A1:N39, matching the refreshed export.
$ws.Rows.Item(9).Delete()
